$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename columns to clean machine-readable names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case the "de"/"de la"/"de las" connector words in a handful of place names
$ws.Range("A14").Value = "Ciudad De México"
$ws.Range("A21").Value = "Estado De México"
$ws.Range("B21").Value = "Naucalpan De Juárez"
$ws.Range("B29").Value = "Silao De La Victoria"
$ws.Range("B32").Value = "Acapulco De Juárez"
$ws.Range("B33").Value = "Chilapa De Álvarez"
$ws.Range("B37").Value = "Pachuca De Soto"
$ws.Range("B39").Value = "Autlán De Navarro"
$ws.Range("B81").Value = "Ignacio De La Llave"

# Drop the trailing metadata/footnote blocks (sample size, source, author, date)
# that lived past the real data range, first the block right after row 91...
$ws.Rows("93:97").Delete()
# ...then the stray duplicate block that used to sit at the very bottom (476:480,
# now shifted up to 471:475 after the first deletion).
$ws.Rows("471:475").Delete()
